# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update "last updated" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 26 de Marzo de 2020 a las 13:12"

# 2. Suiza (row 11) - updated case counts
$ws.Range("B11").Value = 11478
$ws.Range("C11").Value = 581
$ws.Range("E11").Value = 11178
$ws.Range("G11").Value = 16
$ws.Range("H11").Value = 169

# 3. Portugal overtakes Canada and Noruega in ranking -> rows 17-19 reorder
#    Row 17: now Portugal, with updated counts
$ws.Range("A17").Value = "Portugal"
$ws.Range("B17").Value = 3544
$ws.Range("C17").Value = 549
$ws.Range("D17").Value = 43
$ws.Range("E17").Value = 3441
$ws.Range("F17").Value = 61
$ws.Range("G17").Value = 17
$ws.Range("H17").Value = 60

#    Row 18: now Canada, keeping its previous (unchanged) counts
$ws.Range("A18").Value = "Canada"
$ws.Range("B18").Value = 3409
$ws.Range("C18").Value = 0
$ws.Range("D18").Value = 185
$ws.Range("E18").Value = 3188
$ws.Range("F18").Value = 1
$ws.Range("G18").Value = 0
$ws.Range("H18").Value = 36

#    Row 19: now Noruega, with updated counts
$ws.Range("A19").Value = "Noruega"
$ws.Range("B19").Value = 3217
$ws.Range("C19").Value = 133
$ws.Range("D19").Value = 6
$ws.Range("E19").Value = 3197
$ws.Range("F19").Value = 70
$ws.Range("G19").Value = 0
$ws.Range("H19").Value = 14

# 4. Israel (row 23) - updated case counts
$ws.Range("E23").Value = 2423
$ws.Range("G23").Value = 1
$ws.Range("H23").Value = 6

# 5. Pakistan (row 33) - updated case counts
$ws.Range("B33").Value = 1123
$ws.Range("C33").Value = 60
$ws.Range("E33").Value = 1094

# 6. India (row 46) - updated case counts
$ws.Range("B46").Value = 695
$ws.Range("C46").Value = 38
$ws.Range("E46").Value = 636

# 7. Croacia (row 53) - updated case counts
$ws.Range("E53").Value = 457
$ws.Range("G53").Value = 1
$ws.Range("H53").Value = 2

# 8. Hong Kong (row 59) - updated case counts
$ws.Range("F59").Value = 5

# 9. Ucrania (row 85) - updated case counts
$ws.Range("B85").Value = 162
$ws.Range("C85").Value = 17
$ws.Range("E85").Value = 156

# 10. Azerbaiyan (row 93) - updated case counts
$ws.Range("D93").Value = 15
$ws.Range("E93").Value = 104
